# Apply edit: insert a new data row before existing row 220 on the only
# worksheet of the workbook. This pushes old rows 220-223 down to 221-224
# (unchanged), and the newly inserted row 220 receives its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, shifting rows 220:223 down to 221:224.
$ws.Rows.Item(220).Insert()

# New row 220 uses the same layout/style as the rest of the data rows
# (column D uses style index 2 / custom date-time number format). Copy only
# that cell's style from the row below (old row 220, now row 221) so the
# date format on column D carries over correctly without touching the
# whole row / used range.
$ws.Range("D221").Copy()
$ws.Range("D220").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the values for the newly inserted row 220.
$ws.Range("A220").Value = 5
$ws.Range("B220").Value = "Macroferia Regional de Talca"
$ws.Range("C220").Value = "Maule"
$ws.Range("D220").Value = 44595
$ws.Range("E220").Value = 7
$ws.Range("F220").Value = 100112006
$ws.Range("G220").Value = "Repollo"
$ws.Range("H220").Value = "Crespo record"
$ws.Range("I220").Value = "Segunda"
$ws.Range("J220").Value = 2000
$ws.Range("K220").Value = 800
$ws.Range("L220").Value = 800
$ws.Range("M220").Value = 800
$ws.Range("N220").Value = "`$/unidad"
$ws.Range("O220").Value = "Región del Maule"
$ws.Range("P220").Value = 800
$ws.Range("Q220").Value = 1
$ws.Range("R220").Value = "Hortaliza"
